$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function CellRange($table, $row, $col) {
    $c = $table.Cell($row, $col)
    return $d.Range($c.Range.Start, $c.Range.End)
}

function ReplaceInCell($table, $row, $col, $old, $new) {
    $r = CellRange $table $row $col
    $r.Find.ClearFormatting()
    $r.Find.Replacement.ClearFormatting()
    return $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 0, $false, $new, 2)
}

# Row 6 -> version 2.8.1 : ผู้รับผิดชอบ / ผู้ตรวจ
ReplaceInCell $t 6 4 "วรรัตน์ " "ณัฐดนัย" | Out-Null
ReplaceInCell $t 6 4 "(QM)" " (DM)" | Out-Null

ReplaceInCell $t 6 5 "กิตติพศ " "วิรัตน์" | Out-Null
ReplaceInCell $t 6 5 "(SP)" " (TL)" | Out-Null

# Row 7 -> version 1.5.1 : ผู้รับผิดชอบ only
ReplaceInCell $t 7 4 "วรรัตน์ " "วิรัตน์" | Out-Null
ReplaceInCell $t 7 4 "(QM)" " (TL)" | Out-Null

# Row 8 -> version 1.4.2 : ผู้รับผิดชอบ only
ReplaceInCell $t 8 4 "วรรัตน์ " "ณัฐนันท์ " | Out-Null
ReplaceInCell $t 8 4 "(QM)" "(QA)" | Out-Null
